$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAN Data")
$ws.Activate()

# --- Row 26: MPPT poll javed (0x711) - fill in DLC/bytes now known to be all-zero ---
$ws.Range("G26").Value = 8
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Formula = '=IF(ISBLANK(D26),"",1/D26)'

# --- Row 27: MPPT poll woof (0x712) - same update ---
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0

# --- Row 28: Maybe MPPT? / BMS (0x771) - remove old per-byte guesses, note new encoding ---
$ws.Range("H28:O28").Clear()
$ws.Range("G28").Value = 7
$ws.Range("R28").Value = "4*bit&2*bit&2*bit+u_int8+"
$ws.Range("S28").Clear()

# --- Row 29: Maybe MPPT? (0x772) ---
$ws.Range("G29").Value = 7
$ws.Range("R29").ClearContents()
$ws.Range("S29").Clear()

# --- Restore selection / view state ---
$ws.Range("B31").Select()
